# SeniorConnect_MasterLog.xlsx - append newly logged sensor events.
#
# Each worksheet is a simple Date/Timestamp/Hour/Location/Value/Status log.
# We append new rows at the bottom of four of the six sheets (ALERTS,
# mmWave, PIR, Humidity) following the pattern that already exists in
# the sheet. Camera and Proximity are untouched.
#
# All values in this log are stored as literal text (dates like
# "2026-01-30", times like "15:57:22", and percentages like "87.7%" must
# NOT be auto-converted into Excel date/time/number values). We force
# text entry by flipping the cell to a text NumberFormat ("@") right
# before assigning the value, then clearing the formatting again so the
# appended cells keep the same (default) style as the rest of the sheet.

$wb = $excel.ActiveWorkbook

function Append-Rows {
    param(
        [string]$SheetName,
        [int]$StartRow,
        [object[]]$Rows
    )

    $ws = $wb.Worksheets.Item($SheetName)
    $r = $StartRow
    foreach ($rowValues in $Rows) {
        for ($col = 0; $col -lt $rowValues.Length; $col++) {
            $cell = $ws.Cells.Item($r, $col + 1)
            $cell.NumberFormat = "@"
            $cell.Value = $rowValues[$col]
            $cell.ClearFormats()
        }
        $r++
    }
}

# ---------------------------------------------------------------------
# ALERTS: two new FALL_DETECTED critical alerts in the Living Room.
# ---------------------------------------------------------------------
$alertsRows = @(
    @("2026-01-30", "15:57:22", "15:00", "Living Room", "CRITICAL", "FALL_DETECTED"),
    @("2026-01-30", "15:57:52", "15:00", "Living Room", "CRITICAL", "FALL_DETECTED")
)
Append-Rows "ALERTS" 3 $alertsRows

# ---------------------------------------------------------------------
# mmWave: one new presence-detected reading in the Living Room.
# NOTE: a single-row "array of one array" literal gets unwrapped by the
# interpreter unless protected with the leading unary comma, so build it
# that way here (an array of 2+ rows does not need this).
# ---------------------------------------------------------------------
$mmWaveRows = ,@("2026-01-30", "15:57:24", "15:00", "Living Room", "PRESENCE_DETECTED", "Active")
Append-Rows "mmWave" 13 $mmWaveRows

# ---------------------------------------------------------------------
# PIR: a run of bathroom "No Motion" readings with one living-room
# recovery detection mixed in.
# ---------------------------------------------------------------------
$pirRows = @(
    @("2026-01-30", "15:57:11", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:11", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:13", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:18", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:23", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:24", "15:00", "Living Room", "RECOVERY_DETECTION", "Inactive"),
    @("2026-01-30", "15:57:28", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:33", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:38", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:43", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:48", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:53", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:57:58", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:58:03", "15:00", "Bathroom", "No Motion", "Inactive"),
    @("2026-01-30", "15:58:08", "15:00", "Bathroom", "No Motion", "Inactive")
)
Append-Rows "PIR" 33 $pirRows

# ---------------------------------------------------------------------
# Humidity: a run of bathroom humidity percentage readings.
# ---------------------------------------------------------------------
$humidityRows = @(
    @("2026-01-30", "15:57:11", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:57:11", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:57:18", "15:00", "Bathroom", "87.7%", "Active"),
    @("2026-01-30", "15:57:28", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:57:33", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:57:38", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:57:43", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:57:48", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:57:53", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:57:58", "15:00", "Bathroom", "86.2%", "Active"),
    @("2026-01-30", "15:58:03", "15:00", "Bathroom", "87.6%", "Active"),
    @("2026-01-30", "15:58:08", "15:00", "Bathroom", "87.6%", "Active")
)
Append-Rows "Humidity" 28 $humidityRows
